# SectorGroup.xlsx — realign the codeforiati:* columns.
#
# Before:  D = codeforiati:group-name   E = codeforiati:category-name
#          F = codeforiati:category-code G = codeforiati:group-code
# After:   D = codeforiati:category-code E = codeforiati:group-name
#          F = codeforiati:category-name G = codeforiati:group-code
#
# i.e. a 3-way rotation of the D/E/F columns (new D <- old F, new E <- old D,
# new F <- old E) for every row, header included; column G is untouched.
#
# Values are relocated with Range.Copy (cell-to-cell) rather than
# Range.Value/Value2 so that text that looks numeric (sector codes like
# "111", "112", ...) keeps its original string typing instead of being
# coerced into a number, and so no incidental number-format/style is
# introduced on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Range("A1").End(-4121).Row   # xlDown

$colD = "D1:D" + $lastRow
$colE = "E1:E" + $lastRow
$colF = "F1:F" + $lastRow
$scratch = "Z1:Z" + $lastRow

$rangeD = $ws.Range($colD)
$rangeE = $ws.Range($colE)
$rangeF = $ws.Range($colF)
$rangeTemp = $ws.Range($scratch)

# temp = D ; D = F ; F = E ; E = temp
$rangeD.Copy($rangeTemp)
$rangeF.Copy($rangeD)
$rangeE.Copy($rangeF)
$rangeTemp.Copy($rangeE)
$rangeTemp.ClearContents()
